# Update of SInAS version name and of bibliographic citations to data sources
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (AmphRep) keeps the Capinha citation text unchanged; just leave as-is.
# Row 2 (FirstRecords) -> new Seebens citation (Zenodo, v. 3.1)
$ws.Range("N2").Value = "Seebens, H. Alien Species First Records Database. Zenodo https://doi.org/10.5281/ZENODO.10039630 (2023) (v. 3.1)"

# Row 3 (GRIIS) -> new Pagad et al. citation (Zenodo dataset)
$ws.Range("N3").Value = "Shyama Pagad, Bisset, S. & McGeoch, M. A. Country Compendium of the Global Register of Introduced and Invasive Species. Dataset. Zenodo https://doi.org/10.5281/ZENODO.6348164 (2022)"

# Row 5 (GloNAF) -> new Davis et al. citation (Zenodo, v. 3)
$ws.Range("N5").Value = "Davis, A. J. S. et al. Global Naturalized Alien Flora (GloNAF). Open access data to support research on understanding global plant invasions. Zenodo https://doi.org/10.5281/ZENODO.17105725 (2025) (v. 3)"

# Row 6 (GAVIA) -> new Dyer et al. citation (Figshare)
$ws.Range("N6").Value = "Dyer, E., Redding, D. & Blackburn, T. Data from: The Global Avian Invasions Atlas - A database of alien bird distributions worldwide. Figshare https://doi.org/10.6084/M9.FIGSHARE.4234850.V1 (2016)"

# Row 7 (Amph_IUCN) -> IUCN citation updated version string
$ws.Range("N7").Value = "IUCN. The IUCN Red List of Threatened Species. https://www.iucnredlist.org (v. 2024-2)"

# Apply the updated theme-font (Calibri 11, minor scheme) to the refreshed
# citation cells -- matches the new style introduced for these rows.
foreach ($addr in @("N2", "N3", "N5", "N6", "N7")) {
    $ws.Range($addr).Font.ThemeFont = [Microsoft.Office.Interop.Excel.XlThemeFont]::xlThemeFontMinor
}

# Column N is widened to fit the longer citation text; columns L:M keep
# their previous width instead of spanning through N.
$ws.Columns("N:N").ColumnWidth = 93.17

# The active selection moved while editing the citation column.
$ws.Range("N11").Select()
